# Skills Higher Quals fixes.
#
# 1. Shorten the qualification labels in the "Data" sheet's E column from
#    "Diploma" / "Advanced Diploma" to "Dip" / "Adv Dip" (every occurrence,
#    one per each of the six year blocks).
# 2. Widen column E on "Data" now that the new, wider label column needs more
#    room.
# 3. Make "Data" the active sheet/tab again (selecting E25:E26, the last
#    Diploma/Advanced Diploma pair), instead of "Description".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")

# Replace every "Diploma" label with "Dip" ...
$ws1.Range("E5").Value  = "Dip"
$ws1.Range("E9").Value  = "Dip"
$ws1.Range("E13").Value = "Dip"
$ws1.Range("E17").Value = "Dip"
$ws1.Range("E21").Value = "Dip"
$ws1.Range("E25").Value = "Dip"

# ... and every "Advanced Diploma" label with "Adv Dip".
$ws1.Range("E6").Value  = "Adv Dip"
$ws1.Range("E10").Value = "Adv Dip"
$ws1.Range("E14").Value = "Adv Dip"
$ws1.Range("E18").Value = "Adv Dip"
$ws1.Range("E22").Value = "Adv Dip"
$ws1.Range("E26").Value = "Adv Dip"

# Widen column E to fit the relabelled column nicely.
$ws1.Columns.Item(5).ColumnWidth = 37.1

# Switch the active tab/selection back to the Data sheet (it was left on
# Description), selecting the last Diploma/Advanced Diploma row pair.
$ws1.Activate()
$ws1.Range("E25:E26").Select()
